# Insert a new weekly price record for "Bruselas (repollito)" at row 169.
# This pushes the existing rows 169:192 down to 170:193 (values, formats
# and formulas carry down automatically via Insert), then the new row 169
# is populated with the new record's data. Columns that are identical to
# the row that gets pushed down (A,B,C,E,F,G,H,I,O,R) are copied from it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 169:192 down to 170:193, leaving a blank row 169 in place
# (formatting of the row above - e.g. the date style on column D - is
# carried down by Excel's native Insert behaviour).
$ws.Rows.Item(169).Insert()

# Columns that stay the same as the (now shifted) row 170 below.
$ws.Cells.Item(169, 1).Value  = $ws.Cells.Item(170, 1).Value2    # Mercado ID
$ws.Cells.Item(169, 2).Value  = $ws.Cells.Item(170, 2).Value2    # Mercado
$ws.Cells.Item(169, 3).Value  = $ws.Cells.Item(170, 3).Value2    # Región
$ws.Cells.Item(169, 5).Value  = $ws.Cells.Item(170, 5).Value2    # Codreg
$ws.Cells.Item(169, 6).Value  = $ws.Cells.Item(170, 6).Value2    # Categoría ID
$ws.Cells.Item(169, 7).Value  = $ws.Cells.Item(170, 7).Value2    # Categoría
$ws.Cells.Item(169, 8).Value  = $ws.Cells.Item(170, 8).Value2    # Variedad
$ws.Cells.Item(169, 9).Value  = $ws.Cells.Item(170, 9).Value2    # Calidad
$ws.Cells.Item(169, 15).Value = $ws.Cells.Item(170, 15).Value2   # Origen
$ws.Cells.Item(169, 18).Value = $ws.Cells.Item(170, 18).Value2   # Clasificación

# New values specific to this record.
$ws.Cells.Item(169, 4).Value  = 45131                 # Fecha
$ws.Cells.Item(169, 10).Value = 95                    # Volumen
$ws.Cells.Item(169, 11).Value = 24000                 # Precio mínimo
$ws.Cells.Item(169, 12).Value = 24000                 # Precio máximo
$ws.Cells.Item(169, 13).Value = 24000                 # Precio promedio ponderado
$ws.Cells.Item(169, 14).Value = "$/malla 15 kilos"    # Unidad de comercialización
$ws.Cells.Item(169, 16).Value = 1600                  # Precio $/Kg
$ws.Cells.Item(169, 17).Value = 15                    # Kg o Unidades
